# Update membership count values in column B as per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 400
$ws.Range("B3").Value = 500
$ws.Range("B4").Value = 200
$ws.Range("B6").Value = 110
